$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -12.325
$ws.Range("B12").Value = 4.935
$ws.Range("D14").Value = -8.316999999999998
$ws.Range("D19").Value = -7.904999999999999
$ws.Range("C23").Value = -12.247
$ws.Range("D24").Value = -7.449
$ws.Range("B27").Value = 5.246
$ws.Range("C28").Value = -12.604
$ws.Range("B32").Value = 6.407999999999999
$ws.Range("C32").Value = -12.327
$ws.Range("C34").Value = -11.699
$ws.Range("B36").Value = 8.620000000000001
$ws.Range("B38").Value = 5.667
$ws.Range("D38").Value = -7.802
$ws.Range("D41").Value = -8.15
$ws.Range("C42").Value = -12.195
$ws.Range("B46").Value = 6.311
$ws.Range("C49").Value = -12.967
$ws.Range("D52").Value = -7.946
$ws.Range("B54").Value = 5.336
$ws.Range("C54").Value = -12.715
$ws.Range("B55").Value = 4.725999999999999
$ws.Range("B56").Value = 4.632
$ws.Range("B67").Value = 5.502999999999999
$ws.Range("B69").Value = 5.304999999999999
$ws.Range("B72").Value = 5.758
$ws.Range("D72").Value = -7.479000000000001
$ws.Range("C78").Value = -12.156
$ws.Range("D78").Value = -7.382
$ws.Range("C80").Value = -11.332
$ws.Range("B83").Value = 5.113
$ws.Range("D83").Value = -8.608000000000001
$ws.Range("D85").Value = -8.57
$ws.Range("B86").Value = 5.037
$ws.Range("D86").Value = -8.540000000000001
$ws.Range("D90").Value = -6.934
$ws.Range("B91").Value = 5.146
$ws.Range("B93").Value = 5.270000000000001
$ws.Range("D96").Value = -7.229000000000001
$ws.Range("C97").Value = -11.539
$ws.Range("B99").Value = 5.07
$ws.Range("C99").Value = -11.544
$ws.Range("C101").Value = -12.601
$ws.Range("D103").Value = -8.416
$ws.Range("B104").Value = 8.209

$wb.Save()
